$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1 - "I am right on schedule to finish ..." paragraph:
# the phrase "to finish" was duplicated across the run boundary
# ("...to finish " | "to finish the project..."). Drop the duplicated
# leading "to finish " from the second run, and plant the _GoBack bookmark
# at that (relocated) run boundary, right before "the project and present".
# ---------------------------------------------------------------------------

$f1 = $d.Content
$f1.Find.Execute("I am right on schedule to finish ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boundary1 = $f1.End

# Bookmarks.Add re-homes any existing bookmark with the same name, so the
# _GoBack bookmark that currently sits in the "Financially" paragraph moves
# here automatically.
$bmRange1 = $d.Range($boundary1, $boundary1)
$d.Bookmarks.Add("_GoBack", $bmRange1) | Out-Null

# Remove the duplicated "to finish " (10 characters) that starts the run
# right after the new bookmark.
$dupRange = $d.Range($boundary1, $boundary1 + 10)
if ($dupRange.Text -eq "to finish ") {
    $dupRange.Text = ""
}

# ---------------------------------------------------------------------------
# Change 2 - "Financially, ..." paragraph:
# it used to be split in two runs around the _GoBack bookmark
# ("...Therefore, I fo" | bookmark | "cussed solely..."). Now that the
# bookmark has moved away, merge those two runs back into one (pure
# structural merge - no visible text changes). The following run
# ("... and the Raspberry Pi ...") must stay untouched/unmerged, so it is
# temporarily fenced off with a throw-away bookmark while the merge happens.
# ---------------------------------------------------------------------------

$f2 = $d.Content
$f2.Find.Execute("Therefore, I fo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeStart = $f2.Start

$f3 = $d.Content
$f3.Find.Execute("However, the initial budget was for ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$f3.Find.Execute('$200', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$runEnd = $f3.End

# Fence off the boundary that must NOT merge.
$fenceRange = $d.Range($runEnd, $runEnd)
$d.Bookmarks.Add("zzTmpFence", $fenceRange) | Out-Null

# Safety net: if a _GoBack bookmark is still sitting inside the span we are
# about to merge (i.e. Bookmarks.Add above did not already relocate it),
# remove it explicitly.
if ($d.Bookmarks.Exists("_GoBack")) {
    $bmCheck = $d.Bookmarks("_GoBack")
    if ($bmCheck.Start -gt $mergeStart -and $bmCheck.Start -lt $runEnd) {
        $bmCheck.Delete()
    }
}

# Nudge the text so the engine re-normalizes/merges the two adjacent runs
# that now share identical formatting, without altering the visible text:
# insert a sentinel character and immediately delete it again.
$touch = $d.Range($mergeStart, $mergeStart)
$touch.InsertAfter("X")
$sentinel = $d.Range($mergeStart, $mergeStart + 1)
$sentinel.Text = ""

# Remove the temporary fence now that the merge boundary is settled.
$fence = $d.Bookmarks("zzTmpFence")
$fence.Delete()
